$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, pushing existing rows 72..97 down to 73..98.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new record.
$ws.Range("A72").Value = 5
$ws.Range("B72").Value = "Macroferia Regional de Talca"
$ws.Range("C72").Value = "Maule"
$ws.Range("D72").Value = 44489
$ws.Range("E72").Value = 7
$ws.Range("F72").Value = 100112031
$ws.Range("G72").Value = "Poroto verde"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 100
$ws.Range("K72").Value = 43000
$ws.Range("L72").Value = 43000
$ws.Range("M72").Value = 43000
$ws.Range("N72").Value = "$/malla 25 kilos"
$ws.Range("O72").Value = "Región de Arica y Parinacota"
$ws.Range("P72").Value = 1720
$ws.Range("Q72").Value = 25
$ws.Range("R72").Value = "Hortaliza"
